$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right after row 36 (pushing the existing row 37.. data down to 39..68).
$ws.Rows("37:38").Insert()

# Row 37: new data row (2022-03-31, Especial)
$ws.Range("A37").Value = 2
$ws.Range("B37").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C37").Value = "Coquimbo"
$ws.Range("D37").Value = 44651
$ws.Range("D37").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E37").Value = 4
$ws.Range("F37").Value = "Fruta"
$ws.Range("G37").Value = 100107
$ws.Range("H37").Value = "Otros"
$ws.Range("I37").Value = 100107011
$ws.Range("J37").Value = "Tuna"
$ws.Range("K37").Value = "Sin especificar"
$ws.Range("L37").Value = "Especial"
$ws.Range("M37").Value = 240
$ws.Range("N37").Value = 12000
$ws.Range("O37").Value = 13000
$ws.Range("P37").Value = 12500
$ws.Range("Q37").Value = "`$/caja 18 kilos"
$ws.Range("R37").Value = "Provincia de Limarí"
$ws.Range("S37").Value = 694
$ws.Range("T37").Value = 18

# Row 38: new data row (2022-03-31, Primera)
$ws.Range("A38").Value = 2
$ws.Range("B38").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C38").Value = "Coquimbo"
$ws.Range("D38").Value = 44651
$ws.Range("D38").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E38").Value = 4
$ws.Range("F38").Value = "Fruta"
$ws.Range("G38").Value = 100107
$ws.Range("H38").Value = "Otros"
$ws.Range("I38").Value = 100107011
$ws.Range("J38").Value = "Tuna"
$ws.Range("K38").Value = "Sin especificar"
$ws.Range("L38").Value = "Primera"
$ws.Range("M38").Value = 400
$ws.Range("N38").Value = 9000
$ws.Range("O38").Value = 10000
$ws.Range("P38").Value = 9500
$ws.Range("Q38").Value = "`$/caja 18 kilos"
$ws.Range("R38").Value = "Provincia de Limarí"
$ws.Range("S38").Value = 528
$ws.Range("T38").Value = 18

Write-Host "Done. UsedRange rows:"
Write-Host $ws.UsedRange.Rows.Count
